$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row of data (EpiQuery / NY State hospitalization update for 2 April 2020)
$ws.Range("A20").Value = 43924
$ws.Range("B20").Value = 1095
$ws.Range("C20").Value = 395
$ws.Range("D20").Value = 1592
$ws.Range("E20").Value = 351

# Match the date format used by the other cells in column A (copy format only,
# so it reuses the existing style instead of minting a new number format)
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection, matching the saved view state
$ws.Range("D20").Select()
